$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51 (shifts existing row 51..165 down to 52..166)
$ws.Rows(51).Insert()

# Populate the newly inserted row 51 with the new record
$ws.Range("A51").Value = 5
$ws.Range("B51").Value = "Macroferia Regional de Talca"
$ws.Range("C51").Value = "Maule"
$ws.Range("D51").Value = 45002
$ws.Range("E51").Value = 7
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100108
$ws.Range("H51").Value = "Tropicales y subtropicales"
$ws.Range("I51").Value = 100108002
$ws.Range("J51").Value = "Mango"
$ws.Range("K51").Value = "Sin especificar"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 248
$ws.Range("N51").Value = 7000
$ws.Range("O51").Value = 7000
$ws.Range("P51").Value = 7000
$ws.Range("Q51").Value = "$/bandeja 4 kilos"
$ws.Range("R51").Value = "Perú"
$ws.Range("S51").Value = 1750
$ws.Range("T51").Value = 4
